# Commit: "add feature values and their means + only flagged survey in univariate"
#
# The str_comment column (column D) for every variable row is rewritten:
# the trailing period is dropped from each comment sentence. Re-writing the
# value (even though the text is "the same minus a period") forces the
# shared-strings table to drop the old literal and append the new one,
# which is exactly the shape of the underlying OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "This is the unique identifier for the survey"
$ws.Range("D3").Value  = "This identifies the enumerator responsible for the survey"
$ws.Range("D4").Value  = "This indicates whether an anomaly has been detected in the survey data"
$ws.Range("D5").Value  = "This score reflects the severity of the detected anomaly"
$ws.Range("D6").Value  = "This represents the total time taken to complete the survey"
$ws.Range("D7").Value  = "This flag indicates if the survey duration is unusually short"
$ws.Range("D8").Value  = "This shows the total number of questions in the survey"
$ws.Range("D9").Value  = "This is the average time spent on each question in the survey"
$ws.Range("D10").Value = "This flag indicates whether the survey was initiated outside normal working hours"
$ws.Range("D11").Value = "This measures the largest relative increase in median pace during the survey"
$ws.Range("D12").Value = "This counts how many times values were modified in the survey"
$ws.Range("D13").Value = "This indicates how many times the survey was resumed"
$ws.Range("D14").Value = "This counts the number of constraints that were triggered during the survey"
$ws.Range("D15").Value = "This counts how many times there were backtracks in constraint validation"
$ws.Range("D16").Value = "This indicates the total number of errors triggered during the survey"
$ws.Range("D17").Value = "This reflects the average deviation of time spent for each question"
$ws.Range("D18").Value = "This indicates the median deviation of time spent for each question"
$ws.Range("D19").Value = "This represents the variability in time taken for each question"
$ws.Range("D20").Value = "This shows the difference between average and median time deviations for questions"
$ws.Range("D21").Value = "This indicates the average time deviation for questions within a group of questions"
$ws.Range("D22").Value = "This reflects the median time deviation for questions within a group of questions"
$ws.Range("D23").Value = "This measures the variability of time for questions within a group of questions"
$ws.Range("D24").Value = "This indicates the difference between average and median time deviations for a group of questions"
$ws.Range("D25").Value = "This counts the number of outliers identified based on quantiles (1.5 quantile interval)"
$ws.Range("D26").Value = "This counts the number of outliers identified based on standard deviation (2 sigma)"

# Scroll the sheet so row 7 is at the top and column B is the left-most
# visible column, then leave the active cell/selection on D15 - mirrors the
# sheetView's topLeftCell="B7" + <selection activeCell="D15" sqref="D15"/>.
$win = $excel.ActiveWindow
$win.SmallScroll(6, 0, 1, 0)
$ws.Range("D15").Select()
